$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header in F1, matching the style/formatting of the other header cells (B1:E1)
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("F1").Value = "time_taken"

# Add time_taken values for each data row, as text values
$times = @(
    "2021-10-05 13:40:18.378860",
    "2021-10-05 13:40:18.378872",
    "2021-10-05 13:40:18.378876",
    "2021-10-05 13:40:18.378879",
    "2021-10-05 13:40:18.378882",
    "2021-10-05 13:40:18.378885",
    "2021-10-05 13:40:18.378888",
    "2021-10-05 13:40:18.378891"
)

for ($i = 0; $i -lt $times.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $times[$i]
}
